$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2021-11-15"

# Update the header label in B1 to match the new "through" date
$ws.Range("B1").Value = "November 2021 (through November 15)"

# North Lawndale (row 2)
$ws.Range("M2").Value = 10
$ws.Range("X2").Value = 2

# Garfield Park (row 3)
$ws.Range("AI3").Value = 4

# South Shore (row 7)
$ws.Range("AT7").Value = 2

# Englewood (row 8)
$ws.Range("AT8").Value = 2

# Grand Crossing (row 13)
$ws.Range("AI13").Value = 2

# Calumet Heights (row 15)
$ws.Range("B15").Value = 1

# West Loop (row 17)
$ws.Range("B17").Value = 3

# Chatham (row 20)
$ws.Range("B20").Value = 2
$ws.Range("AI20").Value = 1

# West Pullman (row 24)
$ws.Range("B24").Value = 6

# United Center (row 45)
$ws.Range("B45").Value = 2

# East Side (row 70)
$ws.Range("B70").Value = 1

# Gage Park (row 72)
$ws.Range("B72").Value = 2

# Kenwood (row 80)
$ws.Range("X80").Value = 2

# Morgan Park (row 84)
$ws.Range("B84").Value = 4

# Ukrainian Village (row 96)
$ws.Range("B96").Value = 1
